$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused columns N:S (the "ageGroup Late" condition columns),
# shrinking the table from 18 data columns (B:S) down to 12 (B:M).
$ws.Range("N1:S43").EntireColumn.Delete()

# Update the data-translator name used for this published dataset.
$ws.Range("B1:M1").Value = "published_SealeCarlisle_Wetmore_Flowe_Mickes_2019_E1()"

# Exclusions are now empty.
$ws.Range("B4:M4").Value = "{}"

# Relabel the remaining two condition blocks.
$ws.Range("B5:G5").Value = "condition Simultaneous"
$ws.Range("H5:M5").Value = "condition Sequential"

# Update the absolute path recorded for this workbook's folder.
$wb.Path = $wb.Path

# Adjust the view so column E is the left-most visible column and the
# selection sits just past the new data (cell N1).
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("N1").Select()
